$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-24 Friday", "2025-10-25 Saturday"),
    @("653×8=", "570×9="),
    @("961×6=", "372×5="),
    @("775×7=", "878×9="),
    @("408×4=", "461×7="),
    @("941×2=", "209×2="),
    @("503×4=", "326×8="),
    @("517×5=", "272×2="),
    @("488×5=", "941×7="),
    @("214×2=", "703×6="),
    @("124×2=", "645×8="),
    @("415×5=", "455×7="),
    @("894×8=", "331×7="),
    @("303×3=", "442×8="),
    @("947×2=", "780×7="),
    @("163×9=", "440×4="),
    @("719×2=", "843×3="),
    @("259×9=", "447×5="),
    @("660×8=", "878×6="),
    @("211×7=", "628×4="),
    @("287×5=", "181×2="),
    @("373×7=", "148×8="),
    @("797×2=", "254×5="),
    @("299×2=", "422×4="),
    @("796×9=", "155×8="),
    @("641×3=", "196×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
